$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.894.68"
$ws.Range("E2").Value = "  -5.10%  "
$ws.Range("D3").Value = "2.214.96"
$ws.Range("E3").Value = "  -6.41%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.33"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.42%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "99.59"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -8.29%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.590"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -6.41%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.560"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -8.14%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "37.01"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -9.49%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "54.04"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.91%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0829"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -9.35%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.73"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -8.55%  "
$ws.Range("E14").Value = "  -2.45%  "
$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").Value = "2.552.42"
$ws.Range("E15").Value = "  -6.42%  "
$ws.Range("B16").Value = "Polygon"
$ws.Range("C16").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.858"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -12.04%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.18"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -6.81%  "
$ws.Range("D18").Value = "2.211.40"
$ws.Range("E18").Value = "  -6.87%  "
$ws.Range("D19").Value = "42.800.08"
$ws.Range("E19").Value = "  -5.28%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.66"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.31%  "
$ws.Range("D21").Value = "0.0₃0960"
$ws.Range("E21").Value = "  -9.34%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.40"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -10.92%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "65.21"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -10.72%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.13"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -10.72%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "235.50"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -8.85%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.14"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -7.60%  "
$ws.Range("E27").Value = "  -0.13%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.98"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -9.82%  "
$ws.Range("E29").Value = "  -4.72%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.30"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -12.49%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0890"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -8.51%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.51"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -8.20%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "34.38"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -7.42%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "155.58"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -7.52%  "
$ws.Range("E35").Value = "  -7.72%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.14"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +8.63%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.93"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +9.08%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.121"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -6.57%  "
$ws.Range("B39").Value = "NEARProtocol"
$ws.Range("C39").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.91"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.07%  "
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.41"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.61%  "
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.106"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -8.86%  "
$ws.Range("E42").Value = "  -7.78%  "
$ws.Range("D43").Value = "1.909.85"
$ws.Range("E43").Value = "  +0.76%  "
$ws.Range("E44").Value = "  +0.08%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.48"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.01%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "88.34"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -10.88%  "
$ws.Range("E47").Value = "  -9.25%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.37"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.81%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "76.78"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -8.48%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "60.74"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -12.57%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "102.60"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -6.77%  "
